$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4550435.5
$ws.Range("I74").Value = 7147229
$ws.Range("J74").Value = 6047.5
$ws.Range("K74").Value = 7147229
$ws.Range("L74").Value = 6047.5
$ws.Range("M74").Value = -7146293
$ws.Range("N74").Value = -7919.5

$ws.Range("H77").Value = 4550435.5
$ws.Range("I77").Value = 7147229
$ws.Range("J77").Value = 6047.5
$ws.Range("K77").Value = 35736145
$ws.Range("L77").Value = 30237.5
$ws.Range("M77").Value = -35731465
$ws.Range("N77").Value = -39597.5

$ws.Range("H86").Value = 1600
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -77
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 1600
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -384
$ws.Range("N89").Value = -21232

$ws.Range("H138").Value = 3357.6826
$ws.Range("I138").Value = 2340.4443
$ws.Range("J138").Value = 3527.2222
$ws.Range("K138").Value = 7021.3329
$ws.Range("L138").Value = 10581.6666
$ws.Range("M138").Value = -1881.3329
$ws.Range("N138").Value = -20861.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6400.6484
$ws.Range("I32").Value = 4218.3403
$ws.Range("J32").Value = 10199.481
$ws.Range("K32").Value = 4218.3403
$ws.Range("L32").Value = 10199.481
$ws.Range("M32").Value = -3931.3403
$ws.Range("N32").Value = -10773.481

$ws.Range("H122").Value = 1996.0416
$ws.Range("I122").Value = 1444.25
$ws.Range("J122").Value = 4755
$ws.Range("K122").Value = 4332.75
$ws.Range("L122").Value = 14265
$ws.Range("M122").Value = -1882.75
$ws.Range("N122").Value = -19165

$ws.Range("H132").Value = 2594.457
$ws.Range("I132").Value = 1795.3214
$ws.Range("J132").Value = 5791
$ws.Range("K132").Value = 5385.9642
$ws.Range("L132").Value = 17373
$ws.Range("M132").Value = -2855.9642
$ws.Range("N132").Value = -22433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2845.5945
$ws.Range("I134").Value = 1537.8214
$ws.Range("J134").Value = 6914.222
$ws.Range("K134").Value = 4613.4642
$ws.Range("L134").Value = 20742.666
$ws.Range("M134").Value = -2078.4642
$ws.Range("N134").Value = -25812.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1271.7037
$ws.Range("I5").Value = 392.26666
$ws.Range("J5").Value = 2371
$ws.Range("K5").Value = 1176.79998
$ws.Range("L5").Value = 7113
$ws.Range("M5").Value = -1064.79998
$ws.Range("N5").Value = -7337

$ws.Range("H24").Value = 2766.6667
$ws.Range("J24").Value = 4000
$ws.Range("L24").Value = 12000
$ws.Range("N24").Value = -12460

$ws.Range("H118").Value = 2873.625
$ws.Range("I118").Value = 498.16666
$ws.Range("K118").Value = 1494.49998
$ws.Range("M118").Value = -251.4999800000001

$ws.Range("H135").Value = 1271.7037
$ws.Range("I135").Value = 392.26666
$ws.Range("J135").Value = 2371
$ws.Range("K135").Value = 3530.39994
$ws.Range("L135").Value = 21339
$ws.Range("M135").Value = -995.3999400000002
$ws.Range("N135").Value = -26409

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 42475
$ws.Range("J15").Value = 42475
$ws.Range("L15").Value = 42475
$ws.Range("N15").Value = -43051

$ws.Range("H62").Value = 39980
$ws.Range("J62").Value = 39980
$ws.Range("L62").Value = 39980
$ws.Range("N62").Value = -41352

$ws.Range("H64").Value = 25230.066
$ws.Range("J64").Value = 25230.066
$ws.Range("L64").Value = 25230.066
$ws.Range("N64").Value = -25726.066

$ws.Range("H65").Value = 39980
$ws.Range("J65").Value = 39980
$ws.Range("L65").Value = 119940
$ws.Range("N65").Value = -126804

$ws.Range("H67").Value = 25230.066
$ws.Range("J67").Value = 25230.066
$ws.Range("L67").Value = 25230.066
$ws.Range("N67").Value = -26946.066

$ws.Range("H81").Value = 42475
$ws.Range("J81").Value = 42475
$ws.Range("L81").Value = 42475
$ws.Range("N81").Value = -44471

$ws.Range("H84").Value = 42475
$ws.Range("J84").Value = 42475
$ws.Range("L84").Value = 127425
$ws.Range("N84").Value = -137409

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H122").Value = 3022.0476
$ws.Range("I122").Value = 1465.9375
$ws.Range("K122").Value = 4397.8125
$ws.Range("M122").Value = -1947.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15686

$ws.Range("H40").Value = 6561.448
$ws.Range("I40").Value = 4267.9375
$ws.Range("J40").Value = 9384.23
$ws.Range("K40").Value = 4267.9375
$ws.Range("L40").Value = 9384.23
$ws.Range("M40").Value = -4131.9375
$ws.Range("N40").Value = -9656.23

$ws.Range("H122").Value = 4184.9653
$ws.Range("I122").Value = 2760.4285
$ws.Range("J122").Value = 7924.375
$ws.Range("K122").Value = 8281.2855
$ws.Range("L122").Value = 23773.125
$ws.Range("M122").Value = -5831.2855
$ws.Range("N122").Value = -28673.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5103587
$ws.Range("I81").Value = 5496132
$ws.Range("J81").Value = 499
$ws.Range("K81").Value = 10992264
$ws.Range("L81").Value = 998
$ws.Range("M81").Value = -10991203
$ws.Range("N81").Value = -3120

$ws.Range("H84").Value = 5103587
$ws.Range("I84").Value = 5496132
$ws.Range("J84").Value = 499
$ws.Range("K84").Value = 54961320
$ws.Range("L84").Value = 4990
$ws.Range("M84").Value = -54956016
$ws.Range("N84").Value = -15598

$ws.Range("H107").Value = 465.91666
$ws.Range("I107").Value = 509.2
$ws.Range("J107").Value = 249.5
$ws.Range("K107").Value = 1527.6
$ws.Range("L107").Value = 748.5
$ws.Range("M107").Value = 392.4000000000001
$ws.Range("N107").Value = -4588.5

$ws.Range("H122").Value = 3536.1738
$ws.Range("I122").Value = 2109.4285
$ws.Range("J122").Value = 5755.5557
$ws.Range("K122").Value = 6328.2855
$ws.Range("L122").Value = 17266.6671
$ws.Range("M122").Value = -3878.2855
$ws.Range("N122").Value = -22166.6671

$ws.Range("H126").Value = 396594.44
$ws.Range("I126").Value = 1973.9375
$ws.Range("K126").Value = 5921.8125
$ws.Range("M126").Value = -3451.8125

$ws.Range("H132").Value = 13895029
$ws.Range("I132").Value = 7363.8667
$ws.Range("J132").Value = 37041136
$ws.Range("K132").Value = 22091.6001
$ws.Range("L132").Value = 111123408
$ws.Range("M132").Value = -19561.6001
$ws.Range("N132").Value = -111128468

$ws.Range("H135").Value = 67503.42999999999
$ws.Range("J135").Value = 67503.42999999999
$ws.Range("L135").Value = 67503.42999999999
$ws.Range("N135").Value = -77643.42999999999

$ws.Range("H136").Value = 1215.8823
$ws.Range("I136").Value = 474.6154
$ws.Range("J136").Value = 3625
$ws.Range("K136").Value = 1423.8462
$ws.Range("L136").Value = 10875
$ws.Range("M136").Value = 1126.1538
$ws.Range("N136").Value = -15975
